# Auto-generated script applying numeric corrections to Leve profit tables
# across all class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 715.6667
$ws.Range("J17").Value = 703.3043
$ws.Range("L17").Value = 2109.9129
$ws.Range("N17").Value = -2445.9129
$ws.Range("H32").Value = 4775.421
$ws.Range("J32").Value = 1787.4546
$ws.Range("L32").Value = 1787.4546
$ws.Range("N32").Value = -2439.4546
$ws.Range("H33").Value = 320.375
$ws.Range("J33").Value = 455.7
$ws.Range("L33").Value = 455.7
$ws.Range("N33").Value = -913.7
$ws.Range("H116").Value = 19840.834
$ws.Range("I116").Value = 27492.084
$ws.Range("J116").Value = 4538.3335
$ws.Range("K116").Value = 27492.084
$ws.Range("L116").Value = 4538.3335
$ws.Range("M116").Value = -24050.084
$ws.Range("N116").Value = -11422.3335
$ws.Range("H129").Value = 1254.8182
$ws.Range("I129").Value = 793
$ws.Range("K129").Value = 2379
$ws.Range("M129").Value = 2621
$ws.Range("H132").Value = 15906.3125
$ws.Range("I132").Value = 16710.69
$ws.Range("J132").Value = 3840.6667
$ws.Range("K132").Value = 50132.06999999999
$ws.Range("L132").Value = 11522.0001
$ws.Range("M132").Value = -47602.06999999999
$ws.Range("N132").Value = -16582.0001
$ws.Range("H138").Value = 25116.28
$ws.Range("I138").Value = 1449.5883
$ws.Range("J138").Value = 114523.78
$ws.Range("K138").Value = 4348.7649
$ws.Range("L138").Value = 343571.34
$ws.Range("M138").Value = 791.2350999999999
$ws.Range("N138").Value = -353851.34

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15833.721
$ws.Range("I32").Value = 16986.492
$ws.Range("J32").Value = 1308.8
$ws.Range("K32").Value = 16986.492
$ws.Range("L32").Value = 1308.8
$ws.Range("M32").Value = -16699.492
$ws.Range("N32").Value = -1882.8
$ws.Range("H45").Value = 4243.7856
$ws.Range("I45").Value = 2212.1428
$ws.Range("K45").Value = 2212.1428
$ws.Range("M45").Value = -1835.1428
$ws.Range("H61").Value = 8459.066000000001
$ws.Range("I61").Value = 1179.2727
$ws.Range("K61").Value = 1179.2727
$ws.Range("M61").Value = -967.2727
$ws.Range("H122").Value = 2029.35
$ws.Range("I122").Value = 1680.8235
$ws.Range("J122").Value = 4004.3333
$ws.Range("K122").Value = 5042.470499999999
$ws.Range("L122").Value = 12012.9999
$ws.Range("M122").Value = -2592.470499999999
$ws.Range("N122").Value = -16912.9999
$ws.Range("H136").Value = 8459.066000000001
$ws.Range("I136").Value = 1179.2727
$ws.Range("K136").Value = 3537.8181
$ws.Range("M136").Value = -987.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2026.75
$ws.Range("I134").Value = 1363.8125
$ws.Range("J134").Value = 4678.5
$ws.Range("K134").Value = 4091.4375
$ws.Range("L134").Value = 14035.5
$ws.Range("M134").Value = -1556.4375
$ws.Range("N134").Value = -19105.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2349.111
$ws.Range("J94").Value = 3819.7144
$ws.Range("L94").Value = 3819.7144
$ws.Range("N94").Value = -4721.7144
$ws.Range("H99").Value = 5186.25
$ws.Range("I99").Value = 5175.8887
$ws.Range("J99").Value = 5199.5713
$ws.Range("K99").Value = 5175.8887
$ws.Range("L99").Value = 5199.5713
$ws.Range("M99").Value = -3677.8887
$ws.Range("N99").Value = -8195.5713
$ws.Range("H122").Value = 1882.3
$ws.Range("I122").Value = 1721.1111
$ws.Range("K122").Value = 5163.3333
$ws.Range("M122").Value = -2713.3333
$ws.Range("H126").Value = 5186.25
$ws.Range("I126").Value = 5175.8887
$ws.Range("J126").Value = 5199.5713
$ws.Range("K126").Value = 15527.6661
$ws.Range("L126").Value = 15598.7139
$ws.Range("M126").Value = -13057.6661
$ws.Range("N126").Value = -20538.7139
$ws.Range("H132").Value = 60896.65
$ws.Range("I132").Value = 91658.55
$ws.Range("J132").Value = 4499.8335
$ws.Range("K132").Value = 274975.65
$ws.Range("L132").Value = 13499.5005
$ws.Range("M132").Value = -272445.65
$ws.Range("N132").Value = -18559.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 621.0909
$ws.Range("I5").Value = 655.375
$ws.Range("J5").Value = 529.6667
$ws.Range("K5").Value = 1966.125
$ws.Range("L5").Value = 1589.0001
$ws.Range("M5").Value = -1854.125
$ws.Range("N5").Value = -1813.0001
$ws.Range("H33").Value = 565
$ws.Range("I33").Value = 565
$ws.Range("K33").Value = 3390
$ws.Range("M33").Value = -3107
$ws.Range("H131").Value = 100953.4
$ws.Range("J131").Value = 1625.1333
$ws.Range("L131").Value = 4875.3999
$ws.Range("N131").Value = -14955.3999
$ws.Range("H132").Value = 2297.3333
$ws.Range("J132").Value = 2297.3333
$ws.Range("L132").Value = 20675.9997
$ws.Range("N132").Value = -25735.9997
$ws.Range("H135").Value = 621.0909
$ws.Range("I135").Value = 655.375
$ws.Range("J135").Value = 529.6667
$ws.Range("K135").Value = 5898.375
$ws.Range("L135").Value = 4767.0003
$ws.Range("M135").Value = -3363.375
$ws.Range("N135").Value = -9837.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4080.2727
$ws.Range("I122").Value = 4068.1667
$ws.Range("J122").Value = 4094.8
$ws.Range("K122").Value = 12204.5001
$ws.Range("L122").Value = 12284.4
$ws.Range("M122").Value = -9754.500100000001
$ws.Range("N122").Value = -17184.4
$ws.Range("H135").Value = 119990
$ws.Range("J135").Value = 119990
$ws.Range("L135").Value = 119990
$ws.Range("N135").Value = -130130

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1647
$ws.Range("I22").Value = 1472
$ws.Range("J22").Value = 1822
$ws.Range("K22").Value = 1472
$ws.Range("L22").Value = 1822
$ws.Range("M22").Value = -1177
$ws.Range("N22").Value = -2412
$ws.Range("H27").Value = 1647
$ws.Range("I27").Value = 1472
$ws.Range("J27").Value = 1822
$ws.Range("K27").Value = 1472
$ws.Range("L27").Value = 1822
$ws.Range("M27").Value = -1365
$ws.Range("N27").Value = -2036
$ws.Range("H42").Value = 11358.083
$ws.Range("J42").Value = 11358.083
$ws.Range("L42").Value = 11358.083
$ws.Range("N42").Value = -12484.083
$ws.Range("H49").Value = 11358.083
$ws.Range("J49").Value = 11358.083
$ws.Range("L49").Value = 11358.083
$ws.Range("N49").Value = -11652.083
$ws.Range("H82").Value = 1104.2609
$ws.Range("I82").Value = 1223.3572
$ws.Range("J82").Value = 919
$ws.Range("K82").Value = 1223.3572
$ws.Range("L82").Value = 919
$ws.Range("M82").Value = -862.3571999999999
$ws.Range("N82").Value = -1641
$ws.Range("H85").Value = 1104.2609
$ws.Range("I85").Value = 1223.3572
$ws.Range("J85").Value = 919
$ws.Range("K85").Value = 1223.3572
$ws.Range("L85").Value = 919
$ws.Range("M85").Value = 24.64280000000008
$ws.Range("N85").Value = -3415
$ws.Range("H132").Value = 4325
$ws.Range("I132").Value = 3937.5833
$ws.Range("K132").Value = 11812.7499
$ws.Range("M132").Value = -9282.749899999999
$ws.Range("H136").Value = 3111.842
$ws.Range("J136").Value = 3681
$ws.Range("L136").Value = 11043
$ws.Range("N136").Value = -16143

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 10459.8
$ws.Range("J41").Value = 8349.75
$ws.Range("L41").Value = 8349.75
$ws.Range("N41").Value = -9129.75
$ws.Range("H46").Value = 136999
$ws.Range("J46").Value = 136999
$ws.Range("L46").Value = 136999
$ws.Range("N46").Value = -137461
$ws.Range("H132").Value = 5224456.5
$ws.Range("I132").Value = 6115624
$ws.Range("J132").Value = 4758.857
$ws.Range("K132").Value = 18346872
$ws.Range("L132").Value = 14276.571
$ws.Range("M132").Value = -18344342
$ws.Range("N132").Value = -19336.571
$ws.Range("H134").Value = 136999
$ws.Range("J134").Value = 136999
$ws.Range("L134").Value = 410997
$ws.Range("N134").Value = -416067

